$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting rows 13-23 down to 14-24,
# and bringing along correct styles/row-heights for the shifted content.
$ws.Rows.Item(13).Insert()

# The insert leaves a stray empty styled cell at A13; remove it (target row 13 has no A cell).
$ws.Range("A13").Clear()

# New row 13 has no explicit row height, so freshly-created B-column cells default to the
# wrong style; copy the correct style (from B14) into B13 before setting its value.
$ws.Range("B14").Copy($ws.Range("B13"))

# Row 18 should only retain column A (its B/C content moved away); clear the leftovers.
$ws.Range("B18:C18").Clear()

# Fix up the cell contents that were mis-aligned / outdated, per the authored changes.
$ws.Range("B10").Value = 'Apresentar conceitos e metodologias fundamentais para concepção e projeto de fábricas, com ênfase no projeto e organização dos processos de produção, movimentação e armazenagem de materiais, visando o adequado aproveitamento do espaço físico e a eficiência nos fluxos de materiais e ordens de produção no ambiente interno das fábricas.'
$ws.Range("C10").Value = 'Apresentar conceitos e metodologias fundamentais para concepção e projeto de fábricas, com ênfase no projeto e organização dos processos de produção, movimentação e armazenagem de materiais, visando o adequado aproveitamento do espaço físico e a eficiência nos fluxos de materiais e ordens de produção no ambiente interno das fábricas.'
$ws.Range("B13").Value = '5701460 - Antonio Iacono'
$ws.Range("C13").Value = '5701460 - Antonio Iacono'
$ws.Range("B14").Value = 'Arranjo Físico da Fábrica; Sistema de Movimentação e Armazenagem de Materiais; Análise do Fluxo de Materiais.'
$ws.Range("C14").Value = 'Arranjo Físico da Fábrica; Sistema de Movimentação e Armazenagem de Materiais; Análise do Fluxo de Materiais.'
$ws.Range("B16").Value = '1. Objetivos de Desempenho de Empresas de Manufatura e de suas Fábricas. 2. Conceitos de Produto, Recurso e Processo para Projeto da Fábrica. 3. Tipos de Produção e Tipos de Arranjo Físico. 4. Planejamento do Arranjo Físico e dos Fluxos Internos. 5. Manufatura Celular. 6. Planejamento do Sistema de Movimentação e Armazenagem de Materiais'
$ws.Range("C16").Value = '1. Objetivos de Desempenho de Empresas de Manufatura e de suas Fábricas. 2. Conceitos de Produto, Recurso e Processo para Projeto da Fábrica. 3. Tipos de Produção e Tipos de Arranjo Físico. 4. Planejamento do Arranjo Físico e dos Fluxos Internos. 5. Manufatura Celular. 6. Planejamento do Sistema de Movimentação e Armazenagem de Materiais'
$ws.Range("B19").Value = 'Provas e Trabalhos'
$ws.Range("C19").Value = 'Provas e Trabalhos'
$ws.Range("B20").Value = 'M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre projeto de fábricaM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas.A média das provas deve ser maior ou igual a 5,0 (cinco) para que o aluno possa utilizar a nota do Trabalho.'
$ws.Range("C20").Value = 'M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre projeto de fábricaM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas.A média das provas deve ser maior ou igual a 5,0 (cinco) para que o aluno possa utilizar a nota do Trabalho.'
$ws.Range("B21").Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Range("C21").Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Range("B22").Value = 'BANZATO, Eduardo et al. Atualidades na armazenagem. São Paulo: IMAM, 2003.BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. São Paulo, Edgar Blücher, 1977.GURGEL, F.A.C. Administração de recursos materiais e patrimoniais. 2a. Edição. São Paulo. Editora Cengage. 2013. FRANCISCHINI, P.G.; VALLE, C.E. Implantação de Indústrias. Rio de Janeiro, LTC Editora, 1975.LEE, Q et al. Projeto de Instalações e Locais de Trabalho. São Paulo: IMAM, 1998.MOURA, Reinaldo Aparecido. Sistemas e técnicas de movimentação e armazenagem de materiais. IMAM, 2012.NEWMANN, C.; SCALICE, R.K. Projeto de Fábrica e Layout. Rio de Janeiro, Elsevier, 2015.Müther, R. Planejamento do Layout: Sistema SLP. São Paulo, Edgard Blücher, 1978. SLACK, Nigel et al. Administração da produção. São Paulo: Atlas, 8ª ed. 2018.TOMPKINS, James A. et al. Planejamento de instalações. Editora LTC:, 2013.'
$ws.Range("C22").Value = 'BANZATO, Eduardo et al. Atualidades na armazenagem. São Paulo: IMAM, 2003.BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. São Paulo, Edgar Blücher, 1977.GURGEL, F.A.C. Administração de recursos materiais e patrimoniais. 2a. Edição. São Paulo. Editora Cengage. 2013. FRANCISCHINI, P.G.; VALLE, C.E. Implantação de Indústrias. Rio de Janeiro, LTC Editora, 1975.LEE, Q et al. Projeto de Instalações e Locais de Trabalho. São Paulo: IMAM, 1998.MOURA, Reinaldo Aparecido. Sistemas e técnicas de movimentação e armazenagem de materiais. IMAM, 2012.NEWMANN, C.; SCALICE, R.K. Projeto de Fábrica e Layout. Rio de Janeiro, Elsevier, 2015.Müther, R. Planejamento do Layout: Sistema SLP. São Paulo, Edgard Blücher, 1978. SLACK, Nigel et al. Administração da produção. São Paulo: Atlas, 8ª ed. 2018.TOMPKINS, James A. et al. Planejamento de instalações. Editora LTC:, 2013.'

# Narrow the first column definition to column A only (it previously, incorrectly, also
# covered column B). Touching column B's width causes the engine to split the range while
# preserving column A's original exact width.
$ws.Columns.Item(2).ColumnWidth = 60.7109375

